$d = $word.ActiveDocument

# Locate the TRAINING-section paragraph containing the distinctive text
# "Creating API Documentation, LinkedIn Learning Certificate" using Find,
# then resolve which paragraph (by index in $d.Paragraphs) it belongs to.
$findRange = $d.Content
$found = $findRange.Find.Execute("Creating API Documentation", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)

$targetStart = $findRange.Start
$total = $d.Paragraphs.Count
$idx = -1
for ($i = 1; $i -le $total; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if (($targetStart -ge $pp.Range.Start) -and ($targetStart -lt $pp.Range.End)) {
        $idx = $i
        break
    }
}

$p1 = $d.Paragraphs.Item($idx)
$p2 = $d.Paragraphs.Item($idx + 1)

# Change paragraph 1's style from Body Text to List Bullet.
$p1.Style = "List Bullet"

# Change paragraph 2 (the following, currently-empty Body Text paragraph)
# to List Bullet style and give it the new bullet text.
$p2.Style = "List Bullet"
$p2.Range.Text = "Getting Started with GitHub, Pluralsight Training"

$d.Save()
